# Update the "Miguel Andrés-Martínez" speaker textbox on slide 1.
# TextBox 9 (shape id=10) currently reads:
#   "Dr. Miguel Andrés-Martínez "
#   "AWI, Climate Dynamics"
#   "Bremerhaven"
#   <empty>
# It needs to become (name split into separate first/last-name runs,
# "Dr. " title dropped, and "Climate"/rest split into separate runs):
#   "Miguel" " " "Andrés-Martínez" " "
#   "AWI, " "Climate" " Dynamics"
#   "Bremerhaven"
#   <empty>

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 9")

# Keep the shape's on-slide size stable -- it auto-fits its text box
# and we don't want the (legitimate) re-typing below to perturb it.
$origHeight = $shape.Height

$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "Dr. Miguel Andrés-Martínez " -> split name into runs ---
$para1 = $tr.Paragraphs(1)
$para1.Text = "Miguel Andrés-Martínez "
# Re-apply the (unchanged) font size on each sub-range so the single run
# gets split into separate runs at the desired boundaries.
$para1.Characters(1, 6).Font.Size = 24    # "Miguel"
$para1.Characters(7, 1).Font.Size = 24    # " "
$para1.Characters(8, 15).Font.Size = 24   # "Andrés-Martínez"
$para1.Characters(23, 1).Font.Size = 24   # " "

# --- Paragraph 2: "AWI, Climate Dynamics" -> split into 3 runs ---
$para2 = $tr.Paragraphs(2)
$para2.Characters(1, 5).Font.Size = 24    # "AWI, "
$para2.Characters(6, 7).Font.Size = 24    # "Climate"
$para2.Characters(13, 9).Font.Size = 24   # " Dynamics"

# --- Paragraph 3: "Bremerhaven" -> stays a single run ---
$para3 = $tr.Paragraphs(3)
$para3.Characters(1, 11).Font.Size = 24   # "Bremerhaven"

# Restore the shape's original footprint.
$shape.Height = $origHeight
